# Shifting the contents of ShoutOption into functions for portability.
# Rework the options-chain sheet: rename the existing "call"/"put" columns to
# be clearly dated ("Dec 2020 ..."), and add a second expiry ("Jun 2020 ...")
# as two new columns, formatted like the existing price columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "strike"
$ws.Range("B1").Value = "Dec 2020 call"
$ws.Range("C1").Value = "Dec 2020 put"
$ws.Range("D1").Value = "Jun 2020 call"
$ws.Range("E1").Value = "Jun 2020 put"

# --- New data columns: Jun 2020 call / put --------------------------------
$junCall = @(165.1, 161.7, 154.9, 148.3, 141.8, 135.4, 132.2, 129.1, 122.9, 116.9)
$junPut  = @(127.1, 128.6, 131.8, 135.1, 138.5, 142,   143.9, 145.7, 149.5, 153.4)

for ($i = 0; $i -lt $junCall.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 4).Value = $junCall[$i]
    $ws.Cells.Item($row, 5).Value = $junPut[$i]
}

# --- Number formatting for all of the price columns ------------------------
$ws.Range("B2:E11").NumberFormat = "0.00"

# --- Tidy up row heights so they go back to auto-fit ----------------------
$ws.Rows("1:11").AutoFit() | Out-Null

# --- Reset the print/page margins to Excel's normal defaults --------------
$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.7)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.7)
$ws.PageSetup.TopMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.BottomMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.HeaderMargin = $excel.InchesToPoints(0.3)
$ws.PageSetup.FooterMargin = $excel.InchesToPoints(0.3)

# --- Selection moves to G6 as left by the author ---------------------------
$ws.Range("G6").Select() | Out-Null
